# Weekly fruit/vegetable price refresh: the data rows (2-39) get re-keyed to a
# different week. Concretely, the full contents of the "Fecha" (D), "Calidad"
# (I), "Volumen" (J), "Precio minimo" (K), "Precio maximo" (L),
# "Precio promedio ponderado" (M) and "Precio $/Kg" (P) columns move between
# rows according to the mapping below (Target row gets the values that used
# to live in Source row). Columns A,B,C,E,F,G,H,N,O,Q,R are constant across
# every row in this sheet and are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$mapping = @(
    @{Target=2; Source=6},
    @{Target=3; Source=7},
    @{Target=4; Source=29},
    @{Target=5; Source=30},
    @{Target=6; Source=18},
    @{Target=7; Source=37},
    @{Target=8; Source=27},
    @{Target=9; Source=28},
    @{Target=10; Source=14},
    @{Target=11; Source=15},
    @{Target=12; Source=2},
    @{Target=13; Source=3},
    @{Target=14; Source=4},
    @{Target=15; Source=5},
    @{Target=16; Source=10},
    @{Target=17; Source=11},
    @{Target=18; Source=25},
    @{Target=19; Source=26},
    @{Target=20; Source=31},
    @{Target=21; Source=32},
    @{Target=22; Source=12},
    @{Target=23; Source=13},
    @{Target=24; Source=16},
    @{Target=25; Source=17},
    @{Target=26; Source=19},
    @{Target=27; Source=20},
    @{Target=28; Source=23},
    @{Target=29; Source=24},
    @{Target=30; Source=21},
    @{Target=31; Source=22},
    @{Target=32; Source=33},
    @{Target=33; Source=34},
    @{Target=34; Source=38},
    @{Target=35; Source=39},
    @{Target=36; Source=8},
    @{Target=37; Source=9},
    @{Target=38; Source=35},
    @{Target=39; Source=36}
)

# Columns (by index) whose values travel with the row re-keying.
$cols = @(4, 9, 10, 11, 12, 13, 16)   # D, I, J, K, L, M, P

# 1) Snapshot every source value up-front, since a row can be both a
#    "source" for one target and a "target" for another write later on.
$snapshot = @{}
foreach ($col in $cols) {
    for ($row = 2; $row -le 39; $row++) {
        $key = "$row-$col"
        $snapshot[$key] = $ws.Cells.Item($row, $col).Value2
    }
}

# 2) Apply the new values to every target row from the snapshot.
foreach ($m in $mapping) {
    $target = $m.Target
    $source = $m.Source
    foreach ($col in $cols) {
        $key = "$source-$col"
        $ws.Cells.Item($target, $col).Value2 = $snapshot[$key]
    }
}
